$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Text
    if ([string]::IsNullOrEmpty($val)) { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $parts = $val -split ", "

    $idx = -1
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($parts[$i] -ceq "System") {
            $idx = $i
            break
        }
    }
    if ($idx -lt 0) {
        for ($i = 0; $i -lt $parts.Length; $i++) {
            if ($parts[$i] -ceq "admin@admin.com") {
                $idx = $i
                break
            }
        }
    }

    if ($idx -gt 0) {
        $newParts = @($parts[$idx])
        for ($i = 0; $i -lt $parts.Length; $i++) {
            if ($i -ne $idx) { $newParts += $parts[$i] }
        }
        $newVal = $newParts -join ", "
        $cell.Value = $newVal
    }
}
